# Update New Orleans xlsx files:
#  1. Add a "State" column to hotel_info (between Hotel_Name and City) with
#     value "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "State" column into hotel_info ---
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City"; insert a new blank column before it.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder sheet tabs: review_info, then hotel_info ---
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wsHotel)
